$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 67, shifting existing rows 67-73 down to 68-74.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with data (a new price record for Achicoria).
$ws.Range("A67").Value = 10
$ws.Range("B67").Value = "Vega Modelo de Temuco"
$ws.Range("C67").Value = "La Araucanía"
$ws.Range("D67").Value = 45015
$ws.Range("E67").Value = 9
$ws.Range("F67").Value = 100112010
$ws.Range("G67").Value = "Achicoria"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 150
$ws.Range("K67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("M67").Value = 10000
$ws.Range("N67").Value = "$/caja 18 unidades"
$ws.Range("O67").Value = "Región Metropolitana"
$ws.Range("P67").Value = 556
$ws.Range("Q67").Value = 18
$ws.Range("R67").Value = "Hortaliza"
